# Auto-generated edit script: updates profit-calculation columns (H-N)
# on specific rows across multiple sheets, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 408.33334
$ws.Range("I5").Value = 129.16667
$ws.Range("J5").Value = 966.6667
$ws.Range("K5").Value = 129.16667
$ws.Range("L5").Value = 966.6667
$ws.Range("M5").Value = -14.16667000000001
$ws.Range("N5").Value = -1196.6667
$ws.Range("H64").Value = 4444
$ws.Range("I64").Value = 3500.5386
$ws.Range("K64").Value = 3500.5386
$ws.Range("M64").Value = -3252.5386
$ws.Range("H67").Value = 4444
$ws.Range("I67").Value = 3500.5386
$ws.Range("K67").Value = 3500.5386
$ws.Range("M67").Value = -2642.5386
$ws.Range("H74").Value = 5838.926
$ws.Range("J74").Value = 6238.3076
$ws.Range("L74").Value = 6238.3076
$ws.Range("N74").Value = -8110.3076
$ws.Range("H77").Value = 5838.926
$ws.Range("J77").Value = 6238.3076
$ws.Range("L77").Value = 31191.538
$ws.Range("N77").Value = -40551.538
$ws.Range("H125").Value = 3829.7144
$ws.Range("I125").Value = 4586
$ws.Range("K125").Value = 41274
$ws.Range("M125").Value = -38814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1605.4736
$ws.Range("I2").Value = 1244.2307
$ws.Range("J2").Value = 2388.1667
$ws.Range("K2").Value = 1244.2307
$ws.Range("L2").Value = 2388.1667
$ws.Range("M2").Value = -1131.2307
$ws.Range("N2").Value = -2614.1667
$ws.Range("H102").Value = 8400.333000000001
$ws.Range("I102").Value = 8400.333000000001
$ws.Range("K102").Value = 8400.333000000001
$ws.Range("M102").Value = -6778.333000000001
$ws.Range("H116").Value = 1605.4736
$ws.Range("I116").Value = 1244.2307
$ws.Range("J116").Value = 2388.1667
$ws.Range("K116").Value = 1244.2307
$ws.Range("L116").Value = 2388.1667
$ws.Range("M116").Value = 1049.7693
$ws.Range("N116").Value = -6976.1667
$ws.Range("H132").Value = 3449.6182
$ws.Range("I132").Value = 3367.8125
$ws.Range("K132").Value = 10103.4375
$ws.Range("M132").Value = -7573.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1605.4736
$ws.Range("I3").Value = 1244.2307
$ws.Range("J3").Value = 2388.1667
$ws.Range("K3").Value = 1244.2307
$ws.Range("L3").Value = 2388.1667
$ws.Range("M3").Value = -1130.2307
$ws.Range("N3").Value = -2616.1667
$ws.Range("H86").Value = 2043.5
$ws.Range("I86").Value = 1769.4
$ws.Range("J86").Value = 2239.2856
$ws.Range("K86").Value = 1769.4
$ws.Range("L86").Value = 2239.2856
$ws.Range("M86").Value = -646.4000000000001
$ws.Range("N86").Value = -4485.2856
$ws.Range("H89").Value = 2043.5
$ws.Range("I89").Value = 1769.4
$ws.Range("J89").Value = 2239.2856
$ws.Range("K89").Value = 8847
$ws.Range("L89").Value = 11196.428
$ws.Range("M89").Value = -3231
$ws.Range("N89").Value = -22428.428
$ws.Range("H97").Value = 41649.832
$ws.Range("I97").Value = 24857
$ws.Range("K97").Value = 24857
$ws.Range("M97").Value = -23866

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40033.742
$ws.Range("J31").Value = 105334.3
$ws.Range("L31").Value = 105334.3
$ws.Range("N31").Value = -105924.3
$ws.Range("H34").Value = 40033.742
$ws.Range("J34").Value = 105334.3
$ws.Range("L34").Value = 105334.3
$ws.Range("N34").Value = -105738.3
$ws.Range("H86").Value = 6324.8
$ws.Range("J86").Value = 6576
$ws.Range("L86").Value = 6576
$ws.Range("N86").Value = -8822
$ws.Range("H89").Value = 6324.8
$ws.Range("J89").Value = 6576
$ws.Range("L89").Value = 32880
$ws.Range("N89").Value = -44112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 796.35
$ws.Range("I5").Value = 607
$ws.Range("J5").Value = 1553.75
$ws.Range("K5").Value = 1821
$ws.Range("L5").Value = 4661.25
$ws.Range("M5").Value = -1709
$ws.Range("N5").Value = -4885.25
$ws.Range("H33").Value = 1436.5
$ws.Range("I33").Value = 1578.2858
$ws.Range("J33").Value = 444
$ws.Range("K33").Value = 9469.714800000002
$ws.Range("L33").Value = 2664
$ws.Range("M33").Value = -9186.714800000002
$ws.Range("N33").Value = -3230
$ws.Range("H68").Value = 1545.375
$ws.Range("I68").Value = 525
$ws.Range("J68").Value = 1885.5
$ws.Range("K68").Value = 1575
$ws.Range("L68").Value = 5656.5
$ws.Range("M68").Value = -764
$ws.Range("N68").Value = -7278.5
$ws.Range("H71").Value = 1545.375
$ws.Range("I71").Value = 525
$ws.Range("J71").Value = 1885.5
$ws.Range("K71").Value = 4725
$ws.Range("L71").Value = 16969.5
$ws.Range("M71").Value = -669
$ws.Range("N71").Value = -25081.5
$ws.Range("H132").Value = 5931.5
$ws.Range("I132").Value = 3488
$ws.Range("J132").Value = 8375
$ws.Range("K132").Value = 31392
$ws.Range("L132").Value = 75375
$ws.Range("M132").Value = -28862
$ws.Range("N132").Value = -80435
$ws.Range("H135").Value = 796.35
$ws.Range("I135").Value = 607
$ws.Range("J135").Value = 1553.75
$ws.Range("K135").Value = 5463
$ws.Range("L135").Value = 13983.75
$ws.Range("M135").Value = -2928
$ws.Range("N135").Value = -19053.75
$ws.Range("H137").Value = 64465.668
$ws.Range("I137").Value = 88857
$ws.Range("K137").Value = 266571
$ws.Range("M137").Value = -261471

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 209.52942
$ws.Range("I2").Value = 69
$ws.Range("J2").Value = 367.625
$ws.Range("K2").Value = 69
$ws.Range("L2").Value = 367.625
$ws.Range("M2").Value = 44
$ws.Range("N2").Value = -593.625
$ws.Range("H132").Value = 8752.182000000001
$ws.Range("I132").Value = 3899.6667
$ws.Range("J132").Value = 10571.875
$ws.Range("K132").Value = 11699.0001
$ws.Range("L132").Value = 31715.625
$ws.Range("M132").Value = -9169.000100000001
$ws.Range("N132").Value = -36775.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4369.533
$ws.Range("I132").Value = 4208.7334
$ws.Range("J132").Value = 4530.3335
$ws.Range("K132").Value = 12626.2002
$ws.Range("L132").Value = 13591.0005
$ws.Range("M132").Value = -10096.2002
$ws.Range("N132").Value = -18651.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2408.1738
$ws.Range("I132").Value = 2161.3809
$ws.Range("K132").Value = 6484.1427
$ws.Range("M132").Value = -3954.1427
